# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.820.85'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '3.503.52'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''605.82'
$ws.Range('E5').Value = '  +2.60%  '
$ws.Range('D6').Value = '''174.17'
$ws.Range('E6').Value = '  +3.84%  '
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = '3.499.91'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '''0.192'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('E11').Value = '  +7.51%  '
$ws.Range('D12').Value = '''0.582'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '''46.25'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '4.069.24'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '''607.79'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '3.504.94'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '69.835.04'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').Value = '''17.18'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '''0.872'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').Value = '''9.10'
$ws.Range('E23').Value = '  -17.84%  '
$ws.Range('D24').Value = '''15.45'
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').Value = '''95.76'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  -1.27%  '
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').Value = '''34.13'
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').Value = '''8.99'
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').Value = '''684.27'
$ws.Range('E31').Value = '  +19.98%  '
$ws.Range('D32').Value = '''8.12'
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').Value = '''6.93'
$ws.Range('E34').Value = '  +2.01%  '
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').Value = '''0.0999'
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').Value = '''3.56'
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('D38').Value = '''10.70'
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  +8.16%  '
$ws.Range('D40').Value = '''56.47'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('D43').Value = '3.316.62'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').Value = '''2.92'
$ws.Range('E45').Value = '  +4.12%  '
$ws.Range('D46').Value = '''32.28'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('D47').Value = '0.0₃0689'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '''2.56'
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').Value = '''0.130'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').Value = '''133.85'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  -0.08%  '
